$d = $word.ActiveDocument

# --- 1. Remove the stray _GoBack bookmark left after the H-Index paragraph ---
if ($d.Bookmarks.Exists("_GoBack")) {
    $d.Bookmarks("_GoBack").Delete()
}

# --- 2. C-Index paragraph: "increasing" -> "decreasing" ---
$d.Content.Find.Execute(
    "There is a trend of increasing the average C-Index",
    $true, $false, $false, $false, $false, $true, 1, $false,
    "There is a trend of decreasing the average C-Index", 2) | Out-Null

# --- 3. W-Index paragraph: "18 of the 23 stations" -> "19 of the 23 stations" ---
$d.Content.Find.Execute(
    "18 of the 23 stations had positive slopes",
    $true, $false, $false, $false, $false, $true, 1, $false,
    "19 of the 23 stations had positive slopes", 2) | Out-Null

# --- 3b. W-Index paragraph: "1.4 in Winfield" -> "1.6 in Winfield" ---
$d.Content.Find.Execute(
    "1.1 in Saint Francis to 1.4 in Winfield",
    $true, $false, $false, $false, $false, $true, 1, $false,
    "1.1 in Saint Francis to 1.6 in Winfield", 2) | Out-Null

# --- 3c. Re-insert a new _GoBack bookmark right after the new "1.6" ---
$r = $d.Content
$r.Find.Execute("1.6 in Winfield", $true, $false, $false, $false, $false, $true, 1, $false, "", 0) | Out-Null
$bmRange = $d.Range($r.Start + 3, $r.Start + 3)
$d.Bookmarks.Add("_GoBack", $bmRange) | Out-Null
